$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 27.1
$ws.Range("E2").Value = 0.45
$ws.Range("F2").Value = 1627
$ws.Range("I2").Value = 1000
